$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Incadrare" label to "Incadrare juridica" for rows 5-9
$ws.Range("A5:A9").Value = "Incadrare juridica"

# Update numeric values in column B
$ws.Range("B2").Value = 9845
$ws.Range("B3").Value = 2013
$ws.Range("B4").Value = 876
$ws.Range("B5").Value = 432
$ws.Range("B6").Value = 87
$ws.Range("B7").Value = 313
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 28
$ws.Range("B10").Value = 678
$ws.Range("B11").Value = 113
$ws.Range("B12").Value = 45
$ws.Range("B13").Value = 31
$ws.Range("B14").Value = 46
$ws.Range("B15").Value = 371
$ws.Range("B16").Value = 245
$ws.Range("B17").Value = 394
$ws.Range("B18").Value = 87
$ws.Range("B19").Value = 34
$ws.Range("B20").Value = 27
$ws.Range("B21").Value = 14

# Column widths (only column A actually changes width; others keep original)
$ws.Columns.Item(1).ColumnWidth = 17.142857

# Page setup
$ws.PageSetup.PaperSize = 1

# Selection
$ws.Range("B23").Select()
